$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the rows that are no longer part of the table (top-level
#     "Article N - ..." summary rows and the "..." separator rows),
#     keeping only the three rows that already hold the correct
#     Article/Annotation text pairs. Deleting bottom-up keeps the row
#     numbers of earlier deletes stable.
$ws.Rows("7").Delete()
$ws.Rows("5").Delete()
$ws.Rows("1:3").Delete()

# Now the sheet is:
#   Row1: Article 5 - Udemy’s Rights to Content You Post | Udemy should not have the right to use user's content in any way
#   Row2: Article 8.4 Payment and billing                     | Users from the EU should be able to request a refund anytime during a 14 day period
#   Row3: Article 9.3 Limitation of liability                 | Udemy should also be liable for indirect damages and the limit should be uncapped

# --- Insert the new header row ---
$ws.Rows("1").Insert()
$ws.Cells.Item(1,1).Value = "Article"
$ws.Cells.Item(1,2).Value = "Annotation"

# --- Reword the article references (shortened, numbered form) ---
$ws.Cells.Item(2,1).Value = "5. Udemy’s Rights to Content You Post"
$ws.Cells.Item(4,1).Value = "9.3 Limitation of liability"
$ws.Cells.Item(3,1).Value = "8.4 Payments and billing"

# --- Add the extra annotation row for article 8.4 (multiple
#     annotations for the same article), inserted right after the
#     existing 8.4 row. ---
$ws.Rows("4").Insert()
$ws.Cells.Item(4,1).Value = "8.4 Payments and billing"
$ws.Cells.Item(4,2).Value = "Subscription plan can only be changed by user"

# --- Column widths / layout ---
$ws.Columns("A").ColumnWidth = 46.166666666666664
$ws.Columns("B").ColumnWidth = 77

# --- Selection ---
[void]$ws.Range("B2").Select()
